$wb = $excel.ActiveWorkbook

# --- Sheets with a "2015 | 2030 | 2040 | <bug>" header row, where the last
#     header cell ended up holding a stray numeric value instead of the
#     "2050" text label. Fix the label and drop the trailing "Total" row.
$yearLabelSheets = @(
    "Potencia Acumulada - SIN (MW)",
    "Geracao Periodo Medio (MWMed)",
    "Atendimento a Ponta(MW)"
)
foreach ($name in $yearLabelSheets) {
    $ws = $wb.Worksheets.Item($name)

    $e1 = $ws.Range("E1")
    $e1.NumberFormat = "@"
    $e1.Value = "2050"

    $ws.Rows.Item(13).Delete()
}

# --- "Potencia Incremental - SIN(MW)" uses year-range labels
#     (e.g. "2015-2030", "2031-2040"), so its last header becomes "2041-2050".
$ws4 = $wb.Worksheets.Item("Potencia Incremental - SIN(MW)")
$e1_4 = $ws4.Range("E1")
$e1_4.Value = "2041-2050"
$ws4.Rows.Item(13).Delete()

# --- "Emissoes Totais (MtCO2eq)" only needs the header-label fix; it never
#     had a "Total" row to begin with.
$ws5 = $wb.Worksheets.Item("Emissoes Totais (MtCO2eq)")
$e1_5 = $ws5.Range("E1")
$e1_5.NumberFormat = "@"
$e1_5.Value = "2050"

# --- "Custo Total (bilhões de R$)" has no year columns at all, just drop
#     its trailing "Total" row.
$ws6 = $wb.Worksheets.Item("Custo Total (bilhões de R$)")
$ws6.Rows.Item(4).Delete()
